$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "LP & MIP"
